$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2348484848484849
$ws.Range("C2").Value = 0.4962121212121212
$ws.Range("J2").Value = 0.01136363636363636
$ws.Range("P2").Value = 0.1363636363636364
$ws.Range("S2").Value = 0.1212121212121212
$ws.Range("C3").Value = 0.0291970802919708
$ws.Range("J3").Value = 0.0145985401459854
$ws.Range("P3").Value = 0.7299270072992701
$ws.Range("S3").Value = 0.2262773722627737
$ws.Range("J4").Value = 0.03125
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.21875
$ws.Range("B6").Value = 0.05429864253393665
$ws.Range("D6").Value = 0.009049773755656109
$ws.Range("F6").Value = 0.09954751131221719
$ws.Range("J6").Value = 0.1764705882352941
$ws.Range("O6").Value = 0.01809954751131222
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.08144796380090498
$ws.Range("S6").Value = 0.4072398190045249
$ws.Range("B7").Value = 0.1348314606741573
$ws.Range("D7").Value = 0.03370786516853932
$ws.Range("E7").Value = 0.005617977528089887
$ws.Range("F7").Value = 0.03370786516853932
$ws.Range("J7").Value = 0.09550561797752809
$ws.Range("O7").Value = 0.03932584269662921
$ws.Range("Q7").Value = 0.2022471910112359
$ws.Range("R7").Value = 0.101123595505618
$ws.Range("S7").Value = 0.3539325842696629
$ws.Range("B8").Value = 0.06904231625835189
$ws.Range("D8").Value = 0.0111358574610245
$ws.Range("F8").Value = 0.08463251670378619
$ws.Range("J8").Value = 0.1247216035634744
$ws.Range("O8").Value = 0.0155902004454343
$ws.Range("Q8").Value = 0.1937639198218263
$ws.Range("R8").Value = 0.09354120267260579
$ws.Range("S8").Value = 0.4075723830734966
$ws.Range("B9").Value = 0.08823529411764706
$ws.Range("D9").Value = 0.01260504201680672
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.08403361344537816
$ws.Range("O9").Value = 0.02941176470588235
$ws.Range("Q9").Value = 0.2142857142857143
$ws.Range("R9").Value = 0.08403361344537816
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("B10").Value = 0.08753943217665615
$ws.Range("D10").Value = 0.01419558359621451
$ws.Range("F10").Value = 0.06861198738170347
$ws.Range("J10").Value = 0.1198738170347003
$ws.Range("O10").Value = 0.02365930599369085
$ws.Range("Q10").Value = 0.2200315457413249
$ws.Range("R10").Value = 0.09779179810725552
$ws.Range("S10").Value = 0.3682965299684542
$ws.Range("G11").Value = 0.1497975708502024
$ws.Range("J11").Value = 0.08502024291497975
$ws.Range("K11").Value = 0.1740890688259109
$ws.Range("L11").Value = 0.582995951417004
$ws.Range("S11").Value = 0.008097165991902834
$ws.Range("G12").Value = 0.738562091503268
$ws.Range("J12").Value = 0.2091503267973856
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.0261437908496732
$ws.Range("S12").Value = 0.0196078431372549
$ws.Range("G13").Value = 0.7948717948717948
$ws.Range("J13").Value = 0.1794871794871795
$ws.Range("S13").Value = 0.02564102564102564
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02448979591836735
$ws.Range("H15").Value = 0.1591836734693877
$ws.Range("I15").Value = 0.08571428571428572
$ws.Range("J15").Value = 0.363265306122449
$ws.Range("K15").Value = 0.0653061224489796
$ws.Range("M15").Value = 0.00816326530612245
$ws.Range("O15").Value = 0.08571428571428572
$ws.Range("S15").Value = 0.2081632653061224
$ws.Range("F16").Value = 0.02649006622516556
$ws.Range("H16").Value = 0.1986754966887417
$ws.Range("I16").Value = 0.0728476821192053
$ws.Range("J16").Value = 0.3774834437086093
$ws.Range("K16").Value = 0.1258278145695364
$ws.Range("M16").Value = 0.02649006622516556
$ws.Range("O16").Value = 0.06622516556291391
$ws.Range("S16").Value = 0.1059602649006623
$ws.Range("F17").Value = 0.01239669421487603
$ws.Range("H17").Value = 0.1921487603305785
$ws.Range("I17").Value = 0.1053719008264463
$ws.Range("J17").Value = 0.4338842975206612
$ws.Range("K17").Value = 0.08057851239669421
$ws.Range("M17").Value = 0.02272727272727273
$ws.Range("O17").Value = 0.06611570247933884
$ws.Range("S17").Value = 0.08677685950413223
$ws.Range("F18").Value = 0.01785714285714286
$ws.Range("H18").Value = 0.1964285714285714
$ws.Range("I18").Value = 0.1071428571428571
$ws.Range("J18").Value = 0.4375
$ws.Range("K18").Value = 0.08482142857142858
$ws.Range("M18").Value = 0.01785714285714286
$ws.Range("O18").Value = 0.07142857142857142
$ws.Range("S18").Value = 0.06696428571428571
$ws.Range("F19").Value = 0.01082431307243963
$ws.Range("H19").Value = 0.2081598667776852
$ws.Range("I19").Value = 0.1074104912572856
$ws.Range("J19").Value = 0.4004995836802664
$ws.Range("K19").Value = 0.08825978351373855
$ws.Range("M19").Value = 0.01582014987510408
$ws.Range("N19").Value = 0.003330557868442964
$ws.Range("O19").Value = 0.07077435470441298
$ws.Range("S19").Value = 0.09492089925062448
